# "break out stock.yaml completed"
#
# 1) On the "day" sheet, rows 1427..1436 have the bsecode (column D) stored
#    as text; convert those 10 cells to real numbers.
# 2) On the "week" sheet, append 31 new data rows (870..900) for the
#    24/01/2025 weekly pull, and grow the sheet dimension accordingly.
#    The new bsecode values (column D) come in as raw/unprocessed text
#    (not yet cast to numbers) except for row 885 (CDSL), whose bsecode
#    is blank.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fix bsecode column (D) on "day" sheet for rows 1427..1436: text -> number
# ---------------------------------------------------------------------
$dayWs = $wb.Worksheets.Item("day")

$bsecodes = @{
    1427 = 500825
    1428 = 532504
    1429 = 500830
    1430 = 503806
    1431 = 500696
    1432 = 500820
    1433 = 542726
    1434 = 532424
    1435 = 532321
    1436 = 500570
}

foreach ($row in $bsecodes.Keys) {
    $dayWs.Cells.Item($row, 4).Value = $bsecodes[$row]
}

# ---------------------------------------------------------------------
# 2) Append new rows to "week" sheet
# ---------------------------------------------------------------------
$weekWs = $wb.Worksheets.Item("week")

# Each entry: sr, nsecode, name, bsecode, per_chg, close, volume, timeframe, Date Time
$newRows = @(
    @(1,  "PAGEIND",    "Page Industries Limited",                    "532827", -0.96,                 46321.6,   13011,    "week", "24/01/2025 11:32:50"),
    @(2,  "SHREECEM",   "Shree Cements Limited",                      "500387", -0.6899999999999999,   25748,     31358,    "week", "24/01/2025 11:32:50"),
    @(3,  "DIXON",      "Dixon Technologies",                         "540699", 0.78,                   15584.9,   506783,   "week", "24/01/2025 11:32:50"),
    @(4,  "ULTRACEMCO", "Ultratech Cement Limited",                   "532538", -1.18,                  11285.85,  545939,   "week", "24/01/2025 11:32:50"),
    @(5,  "PERSISTENT", "Persistent Systems Limited",                 "533179", 1.16,                   6360.65,   1433020,  "week", "24/01/2025 11:32:50"),
    @(6,  "ABB",        "Abb India Limited",                          "500002", -2.39,                  6178.9,    145838,   "week", "24/01/2025 11:32:50"),
    @(7,  "LTIM",       "LTI Mindtree Ltd",                           "540005", -0.07000000000000001,  5998.15,   219285,   "week", "24/01/2025 11:32:50"),
    @(8,  "POLYCAB",    "Polycab India Ltd",                          "542652", -5.36,                  5915.8,    1354910,  "week", "24/01/2025 11:32:50"),
    @(9,  "SIEMENS",    "Siemens Limited",                            "500550", -3.1,                   5879.7,    247186,   "week", "24/01/2025 11:32:50"),
    @(10, "SUPREMEIND", "Supreme Industries Limited",                 "509930", -4.14,                  3803.2,    176598,   "week", "24/01/2025 11:32:50"),
    @(11, "DMART",      "Avenue Supermarts",                          "540376", -0.51,                  3579.95,   288983,   "week", "24/01/2025 11:32:50"),
    @(12, "TORNTPHARM", "Torrent Pharmaceuticals Limited",            "500420", -0.47,                  3248.4,    391085,   "week", "24/01/2025 11:32:50"),
    @(13, "CUMMINSIND", "Cummins India Limited",                      "500480", -2.94,                  2804.05,   269838,   "week", "24/01/2025 11:32:50"),
    @(14, "NESTLEIND",  "Nestle India Limited",                       "500790", 0.41,                   2207.9,    663600,   "week", "24/01/2025 11:32:50"),
    @(15, "BHARTIARTL", "Bharti Airtel Limited",                      "532454", 0.54,                   1644.8,    3620308,  "week", "24/01/2025 11:32:50"),
    @(16, "CDSL",       "Central Depository Services Ltd",            $null,    -1.05,                  1500.25,   2204447,  "week", "24/01/2025 11:32:50"),
    @(17, "SBILIFE",    "SBI Life Insurance Company Ltd",             "540719", -0.65,                  1440.4,    1504133,  "week", "24/01/2025 11:32:50"),
    @(18, "PRESTIGE",   "Prestige Estates Projects Limited",          "533274", -6.48,                  1255.65,   1730030,  "week", "24/01/2025 11:32:50"),
    @(19, "RELIANCE",   "Reliance Industries Limited",                "500325", -1.37,                  1246.3,    14235970, "week", "24/01/2025 11:32:50"),
    @(20, "ICICIBANK",  "Icici Bank Limited",                         "532174", 0.62,                   1209.2,    9216792,  "week", "24/01/2025 11:32:50"),
    @(21, "LODHA",      "Macrotech Developers Ltd",                   "543287", 1.59,                   1099.55,   2566091,  "week", "24/01/2025 11:32:50"),
    @(22, "TATACONSUM", "TATA Consumer Products Ltd",                 "500800", 0.86,                   992.35,    2201762,  "week", "24/01/2025 11:32:50"),
    @(23, "INDUSINDBK", "Indusind Bank Limited",                      "532187", -2.08,                  950.8,     2938649,  "week", "24/01/2025 11:32:50"),
    @(24, "RAMCOCEM",   "The Ramco Cements Limited",                  "500260", -2.18,                  888.75,    422152,   "week", "24/01/2025 11:32:50"),
    @(25, "PAYTM",      "One 97 Communications Ltd",                  "543396", -4.91,                  807.45,    14201511, "week", "24/01/2025 11:32:50"),
    @(26, "DLF",        "Dlf Limited",                                "532868", -2.77,                  695.25,    3432134,  "week", "24/01/2025 11:32:50"),
    @(27, "JSL",        "Jindal Stainless Limited",                   "532508", -0.88,                  634.1,     597378,   "week", "24/01/2025 11:32:50"),
    @(28, "AUBANK",     "AU Small Finance Bank",                      "540611", 2.1,                    594.65,    5080264,  "week", "24/01/2025 11:32:50"),
    @(29, "VEDL",       "Vedanta Limited",                            "500295", -1.01,                  442,       4758741,  "week", "24/01/2025 11:32:50"),
    @(30, "DELHIVERY",  "Delhivery Ltd",                              "543529", -1.59,                  321.75,    1097240,  "week", "24/01/2025 11:32:50"),
    @(31, "HUDCO",      "Housing and Urban Development Corporation",  "540530", -0.74,                  219.84,    10379337, "week", "24/01/2025 11:32:50")
)

$startRow = 870

# Pre-format column D (bsecode) for the new block as Text, so the
# numeric-looking bsecode values are stored literally as text (matching the
# as-imported state of this freshly appended pull) instead of being
# auto-coerced to numbers. Row 885 (CDSL, blank bsecode) is left out of the
# range so it keeps the plain/default style.
$weekWs.Range("D870:D884").NumberFormat = "@"
$weekWs.Range("D886:D900").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $weekWs.Cells.Item($r, 1).Value = $vals[0]
    $weekWs.Cells.Item($r, 2).Value = $vals[1]
    $weekWs.Cells.Item($r, 3).Value = $vals[2]
    if ($null -ne $vals[3]) {
        $weekWs.Cells.Item($r, 4).Value = $vals[3]
    }
    $weekWs.Cells.Item($r, 5).Value = $vals[4]
    $weekWs.Cells.Item($r, 6).Value = $vals[5]
    $weekWs.Cells.Item($r, 7).Value = $vals[6]
    $weekWs.Cells.Item($r, 8).Value = $vals[7]
    $weekWs.Cells.Item($r, 9).Value = $vals[8]
}
